# Atualizacao de bases das ligas, do dia: 14-06-2024 as 20:31
#
# The underlying data rows (several pairs/triples of rows that share the
# same match Date) had their match-detail columns (B through AD) re-sorted
# among themselves; column A (the running index) is untouched. This script
# reproduces that re-shuffle by swapping/rotating the B:AD content between
# the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peru Liga 1")

function Get-RowValues($sheet, $row) {
    # Columns B (2) through AD (30) as a 1 x 29 variant array.
    return $sheet.Range($sheet.Cells.Item($row, 2), $sheet.Cells.Item($row, 30)).Value2
}

function Set-RowValues($sheet, $row, $values) {
    $sheet.Range($sheet.Cells.Item($row, 2), $sheet.Cells.Item($row, 30)).Value2 = $values
}

# --- Rows 61 and 62: swap content (B:AD) ---
$v61 = Get-RowValues $ws 61
$v62 = Get-RowValues $ws 62
Set-RowValues $ws 61 $v62
Set-RowValues $ws 62 $v61

# --- Rows 175, 176, 177: cyclic rotation ---
# new_175 = old_177, new_176 = old_175, new_177 = old_176
$v175 = Get-RowValues $ws 175
$v176 = Get-RowValues $ws 176
$v177 = Get-RowValues $ws 177
Set-RowValues $ws 175 $v177
Set-RowValues $ws 176 $v175
Set-RowValues $ws 177 $v176

# --- Rows 184, 185, 186: cyclic rotation ---
# new_184 = old_185, new_185 = old_186, new_186 = old_184
$v184 = Get-RowValues $ws 184
$v185 = Get-RowValues $ws 185
$v186 = Get-RowValues $ws 186
Set-RowValues $ws 184 $v185
Set-RowValues $ws 185 $v186
Set-RowValues $ws 186 $v184

# --- Rows 228 and 229: swap content (B:AD) ---
$v228 = Get-RowValues $ws 228
$v229 = Get-RowValues $ws 229
Set-RowValues $ws 228 $v229
Set-RowValues $ws 229 $v228

# --- Rows 252 and 253: swap content (B:AD) ---
$v252 = Get-RowValues $ws 252
$v253 = Get-RowValues $ws 253
Set-RowValues $ws 252 $v253
Set-RowValues $ws 253 $v252
